# "shit ton of fixes"
#
# This particular fix adds a new field/column "c4" (shared-string key "s_i")
# to the p14 sheet:
#   - the existing column Y header "nota_iniciativa" is renamed to "s_i"
#   - a brand new column Z is appended with header "c4"
#   - every data row (2-62) gets a 0 in the new column Z
#
# The new header cell (Z1) should carry the same bold/border/centered
# formatting as the rest of row 1, so we clone it from the neighbouring
# header cell (Y1) rather than re-building the style by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row currently in use (62 in this workbook).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Clone the header formatting from Y1 onto the new Z1 header cell.
$ws.Range("Y1").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rename the Y1 header and set the new Z1 header text.
$ws.Range("Y1").Value = "s_i"
$ws.Range("Z1").Value = "c4"

# Fill the new column with 0 for every existing data row.
$ws.Range("Z2:Z" + $lastRow).Value = 0
